# Update gh-pages to output generated at 456a3b4
# Bumps the "想去人数" (interested-attendee count) figures in column F
# for several rows on the "展览" sheet and the corresponding rows on the
# "全部类型" sheet (which aggregates all event types).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 85
$wsExpo.Range("F3").Value = 812
$wsExpo.Range("F6").Value = 121
$wsExpo.Range("F8").Value = 4619
$wsExpo.Range("F10").Value = 5052
$wsExpo.Range("F11").Value = 577

# --- Sheet "全部类型" (all types combined) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 85
$wsAll.Range("F3").Value = 812
$wsAll.Range("F6").Value = 121
$wsAll.Range("F9").Value = 4619
$wsAll.Range("F11").Value = 5052
$wsAll.Range("F12").Value = 577
